$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# Locate, by content rather than fixed index, the two paragraphs we
# need to touch:
#   - $countParaIdx : the "count = count +1" paragraph currently sitting
#                      BEFORE the "IF max < a[count]" paragraph
#   - $maxParaIdx    : the "max = a[count]" paragraph (inside the IF)
# ------------------------------------------------------------------
$countParaIdx = -1
$maxParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($countParaIdx -eq -1 -and $t -match "count = count \+1") {
        $countParaIdx = $i
    }
    if ($maxParaIdx -eq -1 -and $t -match "max = a\[count\]") {
        $maxParaIdx = $i
    }
}

if ($countParaIdx -eq -1 -or $maxParaIdx -eq -1) {
    Write-Host "Could not locate target paragraphs (count=$countParaIdx max=$maxParaIdx)"
} else {
    # Step 1: remove the "count = count +1" paragraph (with its mark)
    # from its current spot, right before the IF block.
    $d.Paragraphs.Item($countParaIdx).Range.Delete()

    # Re-resolve the "max = a[count]" paragraph index (it shifts down by
    # one once the earlier paragraph is gone).
    if ($maxParaIdx -gt $countParaIdx) {
        $maxParaIdx = $maxParaIdx - 1
    }

    # Step 2: tidy the "max = a[count]" paragraph's run/tab layout -
    # merge the tab immediately preceding the text into that run, and
    # drop the extra trailing tabs.
    $pMax = $d.Paragraphs.Item($maxParaIdx)
    $xmlMax = "<w:p $wns>" +
              "<w:r><w:tab/></w:r>" +
              "<w:r><w:tab/></w:r>" +
              "<w:r><w:tab/><w:t>max = a[count]</w:t></w:r>" +
              "<w:r><w:tab/></w:r>" +
              "</w:p>"
    $pMax.Range.InsertXML($xmlMax)

    # Step 3: insert the "count = count +1" paragraph right after the
    # "max = a[count]" paragraph (i.e. inside the IF block, matching
    # the flowchart image order).
    $pMax = $d.Paragraphs.Item($maxParaIdx)
    $pMax.Range.InsertParagraphAfter()

    $pNew = $d.Paragraphs.Item($maxParaIdx + 1)
    $xmlNew = "<w:p $wns>" +
              "<w:r><w:tab/></w:r>" +
              "<w:r><w:tab/></w:r>" +
              "<w:r><w:t>count = count +1</w:t></w:r>" +
              "<w:r><w:tab/><w:t xml:space=`"preserve`"> </w:t></w:r>" +
              "<w:r><w:tab/></w:r>" +
              "</w:p>"
    $pNew.Range.InsertXML($xmlNew)

    Write-Host "Moved 'count = count +1' paragraph after 'max = a[count]'."
}
